# Update gh-pages output (generated at 456a3b4)
# Refresh "想去人数" (F) counts across sheets, and flip two events that
# became unavailable for sale so their "最低票价" (G) shows "不可售".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 8044
$ws.Range("F3").Value  = 117
$ws.Range("F5").Value  = 30086
$ws.Range("G5").Value  = "不可售"
$ws.Range("F6").Value  = 48
$ws.Range("F7").Value  = 601
$ws.Range("F8").Value  = 684
$ws.Range("F9").Value  = 458
$ws.Range("F10").Value = 144
$ws.Range("F12").Value = 795
$ws.Range("F13").Value = 50
$ws.Range("F14").Value = 597
$ws.Range("F15").Value = 377
$ws.Range("F17").Value = 398
$ws.Range("F18").Value = 149
$ws.Range("F19").Value = 406
$ws.Range("F20").Value = 418
$ws.Range("F21").Value = 1109
$ws.Range("F22").Value = 83
$ws.Range("F23").Value = 692
$ws.Range("F24").Value = 2314
$ws.Range("F25").Value = 816
$ws.Range("F26").Value = 64
$ws.Range("F27").Value = 1079
$ws.Range("F29").Value = 629
$ws.Range("F30").Value = 1070

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value  = 67
$ws.Range("F4").Value  = 336
$ws.Range("F5").Value  = 316
$ws.Range("F10").Value = 2

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 514

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 514
$ws.Range("F3").Value  = 8044
$ws.Range("F4").Value  = 117
$ws.Range("F7").Value  = 30088
$ws.Range("G7").Value  = "不可售"
$ws.Range("F8").Value  = 48
$ws.Range("F9").Value  = 601
$ws.Range("F10").Value = 684
$ws.Range("F11").Value = 458
$ws.Range("F12").Value = 67
$ws.Range("F13").Value = 144
$ws.Range("F15").Value = 336
$ws.Range("F16").Value = 316
$ws.Range("F18").Value = 795
$ws.Range("F19").Value = 50
$ws.Range("F20").Value = 597
$ws.Range("F21").Value = 377
$ws.Range("F26").Value = 2
$ws.Range("F27").Value = 398
$ws.Range("F28").Value = 149
$ws.Range("F29").Value = 406
$ws.Range("F30").Value = 418
$ws.Range("F31").Value = 1109
$ws.Range("F32").Value = 83
$ws.Range("F33").Value = 692
$ws.Range("F34").Value = 2314
$ws.Range("F35").Value = 816
$ws.Range("F36").Value = 64
$ws.Range("F37").Value = 1079
$ws.Range("F40").Value = 629
$ws.Range("F41").Value = 1070
